# "Agenda's tot notulen verbouwen en STD iets aanpassen"
#
# 1) The "STD" (start time) field changes from "11:00" to "1o00" (second
#    "1" retyped as "o", the colon removed) while leaving the text split
#    across three runs ("1" / "o" / "00"), with the document's "_GoBack"
#    bookmark (Word's "last edit location" marker) now sitting between
#    the "o" and the "00".
# 2) Because the "_GoBack" bookmark is a singleton per document, it moves
#    away from its old spot right after "kijken..." near the end of the
#    "Openstaande punten" section.

$d = $word.ActiveDocument

# --- Drop the old "_GoBack" bookmark (after "...kijken...") -------------
# (Word keeps only one "_GoBack" bookmark at a time; it will be re-added
# below at its new location.)
$existing = $d.Bookmarks.Item("_GoBack")
$existing.Delete()

# --- Locate the "11:00" run in the "Tijd" paragraph ----------------------
$timeParagraph = $d.Paragraphs.Item(4)
$timeRange = $timeParagraph.Range
if ($timeRange.Text.Substring(0, 5) -ne "11:00") {
    throw "Expected paragraph 4 to start with '11:00', found '$($timeRange.Text)'"
}
$startPos = $timeRange.Start

# Character layout relative to $startPos:
#   startPos+0 -> "1"   (kept as-is, becomes its own run)
#   startPos+1 -> "1"   (retyped as "o")
#   startPos+2 -> ":"   (deleted)
#   startPos+3..4 -> "00" (kept as-is, becomes its own run)

# Retype the second "1" as "o".
$d.Range($startPos + 1, $startPos + 2).Text = "o"

# Remove the ":" that now sits right after the "o".
$d.Range($startPos + 2, $startPos + 3).Text = ""

# Force a run boundary between "1" and "o" (mirrors the keystroke-by-
# keystroke editing that produced separate runs): drop a throwaway
# bookmark right at that boundary and immediately remove it again -- the
# run split persists even though the bookmark itself is gone.
$splitMark = $d.Bookmarks.Add("zzRunSplit", $d.Range($startPos + 1, $startPos + 1))
$d.Bookmarks.Item("zzRunSplit").Delete()

# Re-add "_GoBack" between "o" and "00", its new home.
$d.Bookmarks.Add("_GoBack", $d.Range($startPos + 2, $startPos + 2)) | Out-Null
